# "now writing a zero to GDX"
# Insert a new row above the existing data on the "index" sheet and add a
# "squeeze=N" label in the newly-widened column F, then make "index" the
# active sheet/tab (it was "sv" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("index")

# Shift the two data rows down by one (old row2 -> row3, old row3 -> row4)
# and leave room for the new row 2.
$ws.Rows.Item(2).Insert()

# New cell: F2 = "squeeze=N"
$ws.Range("F2").Value = "squeeze=N"

# Widen column F to fit the new label (closest reachable width to 18.5703125).
$ws.Columns.Item(6).ColumnWidth = 17.666666666666668

# Make "index" the active sheet/tab and select F2 on it (previously "sv" was
# the active tab with C22 selected).
$ws.Activate() | Out-Null
$ws.Range("F2").Select() | Out-Null
